# Apply the weekly CompStat data refresh described by the commit
# "New crime data collected" to the 112th Precinct workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report volume number + week-covering date range).
# These live inside multi-run shared strings, but every run shares the same
# formatting, so updating the literal text is a faithful, lossless edit.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# ---------------------------------------------------------------------------
# Helper: write a numeric value into a cell (keeps the cell's existing style).
# ---------------------------------------------------------------------------
function Set-Num($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Helper: write a literal text value into a cell that currently holds a
# number, while preserving the "label" style (s=14) used throughout this
# table for non-numeric placeholder cells ("0" / "***.*"). We briefly force
# a text number-format so Excel doesn't re-coerce a numeric-looking string
# back into a number, then repaint the cell's format from a known s=14
# donor cell so the visual style matches the rest of the sheet.
function Set-Label($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-Num "D16" 3
Set-Num "E16" -66.666666666666
Set-Num "F16" 4
Set-Num "G16" 6
Set-Num "H16" -33.333333333333
Set-Num "I16" 37
Set-Num "J16" 40
Set-Num "K16" -7.5
Set-Num "L16" 54.166666666666
Set-Num "M16" -28.846153846153
Set-Num "N16" -89.972899728997

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-Num "C17" 3
Set-Num "D17" 1
Set-Num "E17" 200
Set-Num "F17" 8
Set-Num "H17" 0
Set-Num "I17" 56
Set-Num "J17" 53
Set-Num "K17" 5.660377358490
Set-Num "L17" 64.705882352941
Set-Num "M17" 75
Set-Num "N17" -16.417910447761

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-Num "C18" 1
Set-Num "D18" 2
Set-Num "E18" -50
Set-Num "I18" 66
Set-Num "J18" 57
Set-Num "K18" 15.789473684210
Set-Num "L18" 22.222222222222
Set-Num "M18" -5.714285714285
Set-Num "N18" -90.846047156726

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-Num "C19" 13
Set-Num "D19" 14
Set-Num "E19" -7.142857142857
Set-Num "F19" 45
Set-Num "G19" 33
Set-Num "H19" 36.363636363636
Set-Num "I19" 262
Set-Num "J19" 277
Set-Num "K19" -5.415162454873
Set-Num "L19" 45.555555555555
Set-Num "M19" 40.860215053763
Set-Num "N19" -49.518304431599

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-Num "C20" 6
Set-Num "D20" 2
Set-Num "E20" 200
Set-Num "F20" 21
Set-Num "G20" 11
Set-Num "H20" 90.909090909090
Set-Num "I20" 79
Set-Num "J20" 48
Set-Num "K20" 64.583333333333
Set-Num "L20" 146.875
Set-Num "M20" 46.296296296296
Set-Num "N20" -95.743534482758

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
Set-Num "C21" 24
Set-Num "D21" 22
Set-Num "E21" 9.090909090909
Set-Num "F21" 90
Set-Num "G21" 65
Set-Num "H21" 38.461538461538
Set-Num "I21" 503
Set-Num "J21" 484
Set-Num "K21" 3.925619834710
Set-Num "L21" 51.963746223565
Set-Num "M21" 27.020202020202
Set-Num "N21" -85.782928208027

# ---------------------------------------------------------------------------
# Row 22 - Transit (C22 goes from a numeric 1 to the "0" label placeholder)
# ---------------------------------------------------------------------------
Set-Label "C22" "0"
Set-Num "E22" -100
Set-Num "F22" 4
Set-Num "G22" 3
Set-Num "H22" 33.333333333333
Set-Num "J22" 18
Set-Num "K22" -11.111111111111
Set-Num "M22" 60

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
Set-Num "C24" 30
Set-Num "D24" 29
Set-Num "E24" 3.448275862068
Set-Num "F24" 104
Set-Num "G24" 132
Set-Num "H24" -21.212121212121
Set-Num "I24" 848
Set-Num "J24" 957
Set-Num "K24" -11.389759665621
Set-Num "L24" 28.096676737160
Set-Num "M24" 64.980544747081

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault (D25/E25 become "0"/"***.*" label placeholders)
# ---------------------------------------------------------------------------
Set-Label "D25" "0"
Set-Label "E25" "***.*"
Set-Num "F25" 16
Set-Num "G25" 12
Set-Num "H25" 33.333333333333
Set-Num "I25" 127
Set-Num "K25" 30.927835051546
Set-Num "L25" 30.927835051546
Set-Num "M25" 5.833333333333

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-Num "G26" 1
Set-Num "H26" 100

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes (D27/E27 become "0"/"***.*" label placeholders)
# ---------------------------------------------------------------------------
Set-Label "D27" "0"
Set-Label "E27" "***.*"
Set-Num "F27" 1
Set-Num "H27" -50
